# Updated incorrect screws in parts list
# Applies the cell-level corrections captured by the authoritative diff:
#  - M2.5 screw (row 6): part number / link / price corrections
#  - #2-56 screw (row 9): part number / link / price corrections
#  - Gear motor rows (48/49): add Model/Config (gear ratio) description
#  - Tool rows (106/107/111/112): reorder "Used in Sections" text
#  - Tool rows (108/109/110): point "Used in Sections" at the reworded phrase

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: M2.5 x 4mm Screw -> M2.5 x 6mm Screw, new McMaster part # ---
$ws.Range("A6").Value = "M2.5 x 6mm Screw"
$ws.Range("C6").Value = "92095a458"
$ws.Range("E6").Value = "https://www.mcmaster.com/92095a458"
$ws.Range("I6").Value = 7.07
$ws.Range("J6").Value = 7.07

# --- Row 9: #2-56 x 3/16" Button Head Screws, corrected McMaster part # ---
$ws.Range("C9").Value = "91255A076"
$ws.Range("E9").Value = "https://www.mcmaster.com/91255A076"
$ws.Range("I9").Value = 3.76
$ws.Range("J9").Value = 3.76

# --- Row 48: Gear Motor w Relative Enc. (Drive Motor), add Model/Config ---
$ws.Range("C48").Value = "172:1 gear ratio, with relative encoder"

# --- Row 49: Gear Motor (Corner Motor), add Model/Config ---
$ws.Range("C49").Value = "172:1 gear ratio, NO relative encoder"

# --- Tool rows: reworded "Used in Sections" lists ---
$ws.Range("K106").Value = "Head Assembly, Mechanical Integration, Rocker-Bogie, Differential Pivot, Wheel Assembly, Body"
$ws.Range("K107").Value = "Head Assembly, Mechanical Integration, Rocker-Bogie, Corner Steering, Differential Pivot, Wheel Assembly, Body"
$ws.Range("K108").Value = "Wheel Assembly, Differential Pivot"
$ws.Range("K109").Value = "Wheel Assembly, Differential Pivot"
$ws.Range("K110").Value = "Wheel Assembly, Differential Pivot"
$ws.Range("K111").Value = "Head Assembly, Rocker-Bogie, Corner Steering, Differential Pivot, Wheel Assembly"
$ws.Range("K112").Value = "Corner Steering, Wheel Assembly, Differential Pivot"
